# Apply cryptos list update (commit: Updated cryptos list on Mon Jan 29 09:43:09 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.282.97"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.270.72"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.29"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.21"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.08"
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "2.622.29"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.62"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "2.268.18"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "42.099.50"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.29"
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.64"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.60"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.57"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.98"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.09"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.25"
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0737"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.43"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.115"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("E42").Value = "  -6.73%  "
$ws.Range("D43").Value = "1.958.75"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.98"
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0282"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.83"
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.61"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "2.491.69"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.44"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.51"
$ws.Range("E51").Value = "  -1.80%  "

Write-Host "Applied cryptos update: $([int]89) cell changes"
